$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: "time_taken", styled like the other header cells (copy format from E1)
$ws.Cells.Item(1, 6).Value = "time_taken"
$ws.Cells.Item(1, 5).Copy()
$ws.Cells.Item(1, 6).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data cells F2:F43: timestamp strings (plain, unstyled, forced text)
$ws.Cells.Item(2, 6).Value = "2021-10-05 13:39:15.833105"
$ws.Cells.Item(3, 6).Value = "2021-10-05 13:39:15.833116"
$ws.Cells.Item(4, 6).Value = "2021-10-05 13:39:15.833120"
$ws.Cells.Item(5, 6).Value = "2021-10-05 13:39:15.833123"
$ws.Cells.Item(6, 6).Value = "2021-10-05 13:39:15.833125"
$ws.Cells.Item(7, 6).Value = "2021-10-05 13:39:15.833128"
$ws.Cells.Item(8, 6).Value = "2021-10-05 13:39:15.833131"
$ws.Cells.Item(9, 6).Value = "2021-10-05 13:39:15.833133"
$ws.Cells.Item(10, 6).Value = "2021-10-05 13:39:15.833136"
$ws.Cells.Item(11, 6).Value = "2021-10-05 13:39:15.833139"
$ws.Cells.Item(12, 6).Value = "2021-10-05 13:39:15.833141"
$ws.Cells.Item(13, 6).Value = "2021-10-05 13:39:15.833144"
$ws.Cells.Item(14, 6).Value = "2021-10-05 13:39:15.833146"
$ws.Cells.Item(15, 6).Value = "2021-10-05 13:39:15.833149"
$ws.Cells.Item(16, 6).Value = "2021-10-05 13:39:15.833152"
$ws.Cells.Item(17, 6).Value = "2021-10-05 13:39:15.833155"
$ws.Cells.Item(18, 6).Value = "2021-10-05 13:39:15.833157"
$ws.Cells.Item(19, 6).Value = "2021-10-05 13:39:15.833160"
$ws.Cells.Item(20, 6).Value = "2021-10-05 13:39:15.833163"
$ws.Cells.Item(21, 6).Value = "2021-10-05 13:39:15.833165"
$ws.Cells.Item(22, 6).Value = "2021-10-05 13:39:15.833168"
$ws.Cells.Item(23, 6).Value = "2021-10-05 13:39:15.833170"
$ws.Cells.Item(24, 6).Value = "2021-10-05 13:39:15.833173"
$ws.Cells.Item(25, 6).Value = "2021-10-05 13:39:15.833175"
$ws.Cells.Item(26, 6).Value = "2021-10-05 13:39:15.833178"
$ws.Cells.Item(27, 6).Value = "2021-10-05 13:39:15.833181"
$ws.Cells.Item(28, 6).Value = "2021-10-05 13:39:15.833184"
$ws.Cells.Item(29, 6).Value = "2021-10-05 13:39:15.833186"
$ws.Cells.Item(30, 6).Value = "2021-10-05 13:39:15.833189"
$ws.Cells.Item(31, 6).Value = "2021-10-05 13:39:15.833191"
$ws.Cells.Item(32, 6).Value = "2021-10-05 13:39:15.833194"
$ws.Cells.Item(33, 6).Value = "2021-10-05 13:39:15.833196"
$ws.Cells.Item(34, 6).Value = "2021-10-05 13:39:15.833199"
$ws.Cells.Item(35, 6).Value = "2021-10-05 13:39:15.833202"
$ws.Cells.Item(36, 6).Value = "2021-10-05 13:39:15.833204"
$ws.Cells.Item(37, 6).Value = "2021-10-05 13:39:15.833207"
$ws.Cells.Item(38, 6).Value = "2021-10-05 13:39:15.833209"
$ws.Cells.Item(39, 6).Value = "2021-10-05 13:39:15.833212"
$ws.Cells.Item(40, 6).Value = "2021-10-05 13:39:15.833214"
$ws.Cells.Item(41, 6).Value = "2021-10-05 13:39:15.833217"
$ws.Cells.Item(42, 6).Value = "2021-10-05 13:39:15.833220"
$ws.Cells.Item(43, 6).Value = "2021-10-05 13:39:15.833222"
